$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("D1").Value = "Tiempo estimado"
$ws.Range("E1").Value = "Tiempo real (aproximado)"

# Row 2 - Creación de clases VO (unchanged task/description, add time columns)
$ws.Range("D2").Value = "1 hora"
$ws.Range("E2").Value = "1 hora y 15 minutos"

# Row 3 - Creación de tablas para la BD
$ws.Range("A3").Value = "Creación de tablas para la BD"
$ws.Range("C3").Value = "Creación de fichero con extensión SQL de creación de tablas de la BD."
$ws.Range("D3").Value = "5 minutos"
$ws.Range("E3").Value = "5 minutos"

# Row 4 - NEW row: Pruebas de tablas para la BD
$ws.Range("A4").Value = "Pruebas de tablas para la BD"
$ws.Range("B4").Value = "Terminada"
$ws.Range("C4").Value = "Prueba del fichero con extensión SQL de creación de tablas de la BD."
$ws.Range("D4").Value = "1 minuto"
$ws.Range("E4").Value = "1 minuto"

# Row 5 - Creación de clases DAO (was row 4)
$ws.Range("A5").Value = "Creación de clases DAO"
$ws.Range("B5").Value = "Terminada"
$ws.Range("C5").Value = "Creación de las clases que tendrán funciones que manipulen las tablas de la BD."
$ws.Range("D5").Value = "5 horas"
$ws.Range("E5").Value = "5 horas"

# Row 6 - Creación de clases de prueba I (was row 5, renamed from "Creación de clases de prueba")
$ws.Range("A6").Value = "Creación de clases de prueba I"
$ws.Range("B6").Value = "Terminada"
$ws.Range("C6").Value = "Creación de pruebas para las funciones contenidas en las clases DAO descritas anteriormente."
$ws.Range("D6").Value = "1 hora"
$ws.Range("E6").Value = "1 hora y 30 minutos"

# Row 7 - Creación de una clase fachada (was row 6)
$ws.Range("A7").Value = "Creación de una clase fachada"
$ws.Range("B7").Value = "Terminada"
$ws.Range("C7").Value = "Creación de la clase que contendrá las funciones que usarán los servlets del servidor."
$ws.Range("D7").Value = "5 minutos"
$ws.Range("E7").Value = "5 minutos"

# Row 8 - Creación de servlets (was row 7)
$ws.Range("A8").Value = "Creación de servlets"
$ws.Range("B8").Value = "Terminada"
$ws.Range("C8").Value = "Creación de las clases que tanto la aplicación web como la aplicación Android usarán para proporcionar la funcionalidad especificada de la aplicación."
$ws.Range("D8").Value = "4 horas"
$ws.Range("E8").Value = "4 horas"

# Row 9 - NEW row: Creación de clases de prueba II
$ws.Range("A9").Value = "Creación de clases de prueba II"
$ws.Range("B9").Value = "En progreso"
$ws.Range("C9").Value = "Creación de pruebas para los servlets. "
$ws.Range("D9").Value = "4 horas"

# Row 10 - Instalación de software en servidor (was row 8)
$ws.Range("A10").Value = "Instalación de software en servidor"
$ws.Range("B10").Value = "Terminada"
$ws.Range("C10").Value = "Instalación del SGBD MySQL en el servidor (el entorno Java y el servidor de aplicaciones Apache Tomcat ya estaban instalados)."
$ws.Range("D10").Value = "1 minuto"
$ws.Range("E10").Value = "1 minuto"

# Row 11 - Configuración de software en servidor (was row 9)
$ws.Range("A11").Value = "Configuración de software en servidor"
$ws.Range("B11").Value = "En progreso"
$ws.Range("C11").Value = "Configuración del SGBD MySQL (principalmente usuarios admitidos) y del servidor de aplicaciones Apache Tomcat."
$ws.Range("D11").Value = "30 minutos"

# Row 12 - NEW row: Adaptación de los servlets al formato JSON
$ws.Range("A12").Value = "Adaptación de los servlets al formato JSON"
$ws.Range("B12").Value = "En progreso"
$ws.Range("C12").Value = "Los servlets se adaptarán para que puedan servir tanto para la aplicación web como para la aplicación Android."
$ws.Range("D12").Value = "1 hora"

# Column widths to match final layout (values chosen so the stored OOXML
# "width" attribute - which adds ~5/6 character of padding on top of the
# COM ColumnWidth - lands on the target widths used by the final workbook).
$ws.Columns.Item(1).ColumnWidth = 36
$ws.Columns.Item(3).ColumnWidth = 125.33333333333333
$ws.Columns.Item(4).ColumnWidth = 18.5
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# Apply header style (bold) to the new header cells D1:E1, reusing the
# existing header style (same as A1:C1) instead of creating a new one.
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection to match final state
$ws.Range("D12").Select()
